$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.704.36'
$ws.Range('E2').Value = '  -5.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.545.11'
$ws.Range('E3').Value = '  -2.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '299.77'
$ws.Range('E5').Value = '  -2.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.56'
$ws.Range('E6').Value = '  -5.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.578'
$ws.Range('E7').Value = '  -3.97%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.559'
$ws.Range('E9').Value = '  -3.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.91'
$ws.Range('E10').Value = '  -6.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0813'
$ws.Range('E11').Value = '  -4.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.76'
$ws.Range('E12').Value = '  -5.26%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.107'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.932.33'
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.568.42'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.885'
$ws.Range('E16').Value = '  -4.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.26'
$ws.Range('E17').Value = '  -5.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.676.45'
$ws.Range('E18').Value = '  -6.17%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.65'
$ws.Range('E19').Value = '  -1.72%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0978'
$ws.Range('E20').Value = '  -3.70%  '
$ws.Range('E21').Value = '  -3.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.40'
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '263.84'
$ws.Range('E23').Value = '  -4.02%  '
$ws.Range('E24').Value = '  -3.82%  '
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.18'
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.20'
$ws.Range('E28').Value = '  -4.29%  '
$ws.Range('E29').Value = '  -2.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.77'
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.13'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('E32').Value = '  -4.05%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.79'
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '151.44'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.17'
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0808'
$ws.Range('E36').Value = '  -4.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  -4.50%  '
$ws.Range('E38').Value = '  -2.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.60'
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.36'
$ws.Range('E40').Value = '  +2.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.56'
$ws.Range('E41').Value = '  -2.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0314'
$ws.Range('E42').Value = '  -5.73%  '
$ws.Range('E43').Value = '  -5.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.023.91'
$ws.Range('E44').Value = '  -5.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '87.60'
$ws.Range('E46').Value = '  -6.25%  '
$ws.Range('E47').Value = '  +5.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.07'
$ws.Range('E48').Value = '  -4.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.795.44'
$ws.Range('E49').Value = '  -2.49%  '
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('E51').Value = '  -5.20%  '
